$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row -> new Price (column D) value, as text (matches original inline-string formatting)
$updates = @{
    2  = "261.17"
    3  = "22.85"
    4  = "6.204"
    5  = "0.06140"
    6  = "6.733"
    7  = "3.463"
    8  = "1.371"
    9  = "0.7991"
    10 = "0.1580"
    11 = "0.08110"
    12 = "0.03472"
    13 = "0.03069"
    14 = "0.09316"
    15 = "3.847"
    16 = "0.001687"
    18 = "0.0006150"
    19 = "0.006213"
    21 = "0.004061"
    23 = "3.694"
    24 = "2.216"
    40 = "0.04611"
    41 = "0.007092"
    43 = "0.003131"
    44 = "0.01024"
    46 = "0.00005942"
    48 = "0.7000"
    49 = "0.08983"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$row]
    $cell.Style = "Normal"
}
